# Update input cost data on the "time_variants" sheet (row 3 = program_cost_vaccination)
# and reflect the resulting UI state (active sheet/tab, selection, column width).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")

# --- Updated input cost values (row 3) ---
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 20000
$ws.Range("O3").Value = 30000
$ws.Range("T3").Value = 50000
$ws.Range("V3").Value = 70000
$ws.Range("X3").Value = 100000
$ws.Range("Y3").Value = 120000
$ws.Range("Z3").Value = 100000
$ws.Range("AA3").Value = 100000
$ws.Range("AB3").Value = 120000
$ws.Range("AC3").Value = 150000
$ws.Range("AD3").Value = 200000
$ws.Range("AE3").Value = 210000
$ws.Range("AF3").Value = 250000
$ws.Range("AG3").Value = 250000
$ws.Range("AH3").Value = 250000

# Column AE widened to fit the new best-fit value
$ws.Columns.Item(31).AutoFit()

# Make "time_variants" the active/selected sheet (was "constants"),
# and move the selection on it from G8 to A6.
$ws.Activate()
[void]$ws.Range("A6").Select()
